$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.210.46'
$ws.Range('E2').Value = '  -1.47%  '

$ws.Range('D3').Value = '3.785.60'
$ws.Range('E3').Value = '  +0.11%  '

$ws.Range('E4').Value = '  -1.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '419.02'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.79'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -7.02%  '

$ws.Range('D7').Value = '3.786.34'
$ws.Range('E7').Value = '  -3.74%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.68%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.715'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.160'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -9.91%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000342'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -6.34%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.59'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -6.55%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.371.02'
$ws.Range('E14').Value = '  -0.63%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.99'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.19%  '

$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.93'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +19.78%  '

$ws.Range('E17').Value = '  -1.70%  '

$ws.Range('D18').Value = '3.786.83'
$ws.Range('E18').Value = '  +1.40%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.48%  '

$ws.Range('D20').Value = '66.295.41'
$ws.Range('E20').Value = '  -1.65%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.07'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '403.62'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -7.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.10'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.94'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -6.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.99%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '36.61'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.98%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.71'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +10.45%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.50%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.27'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '696.52'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.02'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +11.80%  '

$ws.Range('E32').Value = '  +1.04%  '

$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.17'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.72%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.119'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.149'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -7.21%  '

$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '37.30'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -9.69%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '54.45'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.66%  '

$ws.Range('D39').Value = '0.0₃0758'
$ws.Range('E39').Value = '  +13.40%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0448'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.06%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.09%  '

$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.69'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +8.35%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('E44').Value = '  -7.75%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.30'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.07%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '143.82'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.38%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.70%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.02'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.10%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.30'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.28%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.72'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.71%  '

$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.49'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.17%  '
